$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price (column D) values that look numeric keep their original
# textual formatting (leading/trailing zeros, multi-dot separators, etc.)
# by forcing the cell to Text format before assigning the string value.

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '25.095.11'
$ws.Range('E2').Value = '  -2.87%  '

$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.654.74'
$ws.Range('E3').Value = '  -4.67%  '

$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.000'
$ws.Range('E4').Value = '  -0.03%  '

$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '237.00'
$ws.Range('E5').Value = '  -1.68%  '

$ws.Range('E6').Value = '  -0.03%  '

$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.4793'
$ws.Range('E7').Value = '  -7.49%  '

$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.2628'
$ws.Range('E8').Value = '  -4.03%  '

$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.05986'
$ws.Range('E9').Value = '  -2.61%  '

$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.07109'

$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '1.657.38'
$ws.Range('E11').Value = '  -4.61%  '

$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '14.45'
$ws.Range('E12').Value = '  -3.12%  '

$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.6205'
$ws.Range('E13').Value = '  -3.01%  '

$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '4.594'
$ws.Range('E14').Value = '  -0.18%  '

$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '73.13'
$ws.Range('E15').Value = '  -5.19%  '

$ws.Range('E16').Value = '  -0.02%  '

$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '0.9996'
$ws.Range('E17').Value = '  -0.07%  '

$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '25.078.73'
$ws.Range('E18').Value = '  -3.07%  '

$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '11.38'
$ws.Range('E19').Value = '  -2.69%  '

$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000006566'
$ws.Range('E20').Value = '  -2.88%  '

$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '4.454'
$ws.Range('E21').Value = '  +4.61%  '

$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '1.865.93'
$ws.Range('E22').Value = '  -4.97%  '

$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '8.469'
$ws.Range('E23').Value = '  -1.79%  '

$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '5.263'
$ws.Range('E24').Value = '  +0.30%  '

$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '133.11'
$ws.Range('E25').Value = '  -3.57%  '

$ws.Range('E26').Value = '  -3.02%  '

$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '1.392'
$ws.Range('E27').Value = '  -8.36%  '

$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '1.693'
$ws.Range('E28').Value = '  -4.15%  '

$ws.Range('E29').Value = '  -3.16%  '

$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.804'
$ws.Range('E30').Value = '  -3.56%  '

$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '0.07913'
$ws.Range('E31').Value = '  -4.08%  '

$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.525'
$ws.Range('E32').Value = '  -3.49%  '

$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.04593'
$ws.Range('E33').Value = '  -0.74%  '

$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.603'
$ws.Range('E34').Value = '  -1.54%  '

$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '0.9399'
$ws.Range('E35').Value = '  -4.56%  '

$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.5818'
$ws.Range('E36').Value = '  -5.47%  '

$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.624'
$ws.Range('E37').Value = '  -2.26%  '

$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01537'
$ws.Range('E38').Value = '  -3.79%  '

$ws.Range('E39').Value = '  +12.76%  '

$ws.Range('E40').Value = '  +0.06%  '

$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '1.827'
$ws.Range('E41').Value = '  -4.61%  '

$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '98.62'
$ws.Range('E42').Value = '  -1.22%  '

$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '0.3704'
$ws.Range('E43').Value = '  -3.40%  '

$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '4.831'
$ws.Range('E44').Value = '  -3.20%  '

$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.1118'
$ws.Range('E45').Value = '  -0.41%  '

$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '6.060'
$ws.Range('E46').Value = '  -2.72%  '

$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '0.05144'
$ws.Range('E47').Value = '  -1.88%  '

$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '52.13'
$ws.Range('E48').Value = '  -4.74%  '

$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '29.35'
$ws.Range('E49').Value = '  -3.74%  '

$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').Value = '  -0.10%  '

$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.9991'
$ws.Range('E51').Value = '  +0.07%  '
